# Adds rows 190-224 to the "system logs" worksheet, recording the new
# CPF (copy file) and MVF (move file) command log entries described in the
# commit "adding CPF & MVF commands."

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Template cells already carrying the workbook's conditional-style colors:
#   F2 -> green "Success" style,  F3 -> red "Failed" style
$successTemplate = $ws.Range("F2")
$failTemplate = $ws.Range("F3")

$logRows = @(
    @{ A = "2021/05/16"; B = "01:54:41"; C = "mjavadtatari"; D = "login"; E = "pass: mjavadtatari"; F = "Success"; Style = 2 },
    @{ A = "2021/05/16"; B = "01:54:46"; C = "mjavadtatari"; D = "DEF"; E = "root\mjavadtatari\jj.txt"; F = "Failed. File Not Found!"; Style = 3 },
    @{ A = "2021/05/16"; B = "01:54:54"; C = "mjavadtatari"; D = "HOME"; E = "root\mjavadtatari\home\"; F = "Directory Changed Successfully"; Style = 2 },
    @{ A = "2021/05/16"; B = "01:54:57"; C = "mjavadtatari"; D = "DEF"; E = "root\mjavadtatari\home\jj.txt"; F = "Success. File Removed!"; Style = 2 },
    @{ A = "2021/05/16"; B = "01:55:00"; C = "mjavadtatari"; D = "DEF"; E = "root\mjavadtatari\home\jj.txt"; F = "Failed. File Not Found!"; Style = 3 },
    @{ A = "2021/05/16"; B = "01:56:09"; C = "mjavadtatari"; D = "login"; E = "pass: mjavadtatari"; F = "Success"; Style = 2 },
    @{ A = "2021/05/16"; B = "01:56:29"; C = "mjavadtatari"; D = "DEF"; E = "root\U1030\home\d.txt"; F = "Success. File Removed!"; Style = 2 },
    @{ A = "2021/05/16"; B = "01:58:01"; C = "mjavadtatari"; D = "login"; E = "pass: mjavadtatari"; F = "Success"; Style = 2 },
    @{ A = "2021/05/16"; B = "01:58:11"; C = "mjavadtatari"; D = "DEF"; E = "root\mjavadtatari\U1030\home\asd.txt"; F = "Failed. File Not Found!"; Style = 3 },
    @{ A = "2021/05/16"; B = "02:00:20"; C = "mjavadtatari"; D = "login"; E = "pass: mjavadtatari"; F = "Success"; Style = 2 },
    @{ A = "2021/05/16"; B = "02:00:29"; C = "mjavadtatari"; D = "MKD"; E = "U1030\home\hell\"; F = "Fail. Not Valid!"; Style = 3 },
    @{ A = "2021/05/16"; B = "02:00:54"; C = "mjavadtatari"; D = "DED"; E = "root\mjavadtatari\U1030\home\"; F = "Failed. Directory Does Not Exists!"; Style = 3 },
    @{ A = "2021/05/16"; B = "02:02:53"; C = "mjavadtatar"; D = "login"; E = "pass: mjavadtatar"; F = "Failed Attempts= 1"; Style = 3 },
    @{ A = "2021/05/16"; B = "02:02:55"; C = "mjavadtatari"; D = "login"; E = "pass: mjavadtatari"; F = "Success"; Style = 2 },
    @{ A = "2021/05/16"; B = "02:03:03"; C = "mjavadtatari"; D = "HELP"; E = "help cpd"; F = "Success"; Style = 2 },
    @{ A = "2021/05/16"; B = "02:09:51"; C = "mjavadtatari"; D = "login"; E = "pass: mjavadtatari"; F = "Success"; Style = 2 },
    @{ A = "2021/05/16"; B = "02:10:06"; C = "mjavadtatari"; D = "CPF"; E = ", "; F = "Failed. Invalid Parameters!"; Style = 3 },
    @{ A = "2021/05/16"; B = "02:11:34"; C = "mjavadtatari"; D = "login"; E = "pass: mjavadtatari"; F = "Success"; Style = 2 },
    @{ A = "2021/05/16"; B = "02:11:37"; C = "mjavadtatari"; D = "HOME"; E = "root\mjavadtatari\home\"; F = "Directory Changed Successfully"; Style = 2 },
    @{ A = "2021/05/16"; B = "02:11:49"; C = "mjavadtatari"; D = "CPF"; E = "root\mjavadtatari\home\1.txt --> root\mjavadtatari\home\2.txt"; F = "Success. File Copied!"; Style = 2 },
    @{ A = "2021/05/16"; B = "02:12:01"; C = "mjavadtatari"; D = "CPF"; E = "root\mjavadtatari\home\1.txt --> root\mjavadtatari\home\2.txt"; F = "Success. File Copied!"; Style = 2 },
    @{ A = "2021/05/16"; B = "02:13:31"; C = "mjavadtatari"; D = "login"; E = "pass: mjavadtatari"; F = "Success"; Style = 2 },
    @{ A = "2021/05/16"; B = "02:13:51"; C = "mjavadtatari"; D = "CPF"; E = "root\mjavadtatari\home\1.txt --> root\mjavadtatari\home\3.txt"; F = "Success. File Copied!"; Style = 2 },
    @{ A = "2021/05/16"; B = "02:13:54"; C = "mjavadtatari"; D = "CPF"; E = "root\mjavadtatari\home\1.txt --> root\mjavadtatari\home\3.txt"; F = "Failed. File Already Exists!"; Style = 3 },
    @{ A = "2021/05/16"; B = "02:15:29"; C = "mjavadtatari"; D = "login"; E = "pass: mjavadtatari"; F = "Success"; Style = 2 },
    @{ A = "2021/05/16"; B = "02:15:46"; C = "mjavadtatari"; D = "CPF"; E = "root\mjavadtatari\home\12.txt --> root\mjavadtatari\home\5.txt"; F = "Failed. File Not Found!"; Style = 3 },
    @{ A = "2021/05/16"; B = "02:15:51"; C = "mjavadtatari"; D = "CPF"; E = "root\mjavadtatari\home\1.txt --> root\mjavadtatari\home\5.txt"; F = "Success. File Copied!"; Style = 2 },
    @{ A = "2021/05/16"; B = "02:15:53"; C = "mjavadtatari"; D = "CPF"; E = "root\mjavadtatari\home\1.txt --> root\mjavadtatari\home\5.txt"; F = "Failed. File Already Exists!"; Style = 3 },
    @{ A = "2021/05/16"; B = "02:18:20"; C = "mjavadtatari"; D = "login"; E = "pass: mjavadtatari"; F = "Success"; Style = 2 },
    @{ A = "2021/05/16"; B = "02:18:22"; C = "mjavadtatari"; D = "HOME"; E = "root\mjavadtatari\home\"; F = "Directory Changed Successfully"; Style = 2 },
    @{ A = "2021/05/16"; B = "02:18:23"; C = "mjavadtatari"; D = "SHOW"; E = "show"; F = "Success. Sub-Directories and Files Showed-Up!"; Style = 2 },
    @{ A = "2021/05/16"; B = "02:18:45"; C = "mjavadtatari"; D = "MVF"; E = "root\mjavadtatari\home\3.txt --> root\mjavadtatari\home\well\3.txt"; F = "Success. File Moved!"; Style = 2 },
    @{ A = "2021/05/16"; B = "02:19:00"; C = "mjavadtatari"; D = "SHOW"; E = "show"; F = "Success. Sub-Directories and Files Showed-Up!"; Style = 2 },
    @{ A = "2021/05/16"; B = "02:19:02"; C = "mjavadtatari"; D = "MVF"; E = "root\mjavadtatari\home\3.txt --> root\mjavadtatari\home\well\3.txt"; F = "Failed. File Not Found!"; Style = 3 },
    @{ A = "2021/05/16"; B = "02:19:13"; C = "mjavadtatari"; D = "MVF"; E = "root\mjavadtatari\home\2.txt --> root\mjavadtatari\home\well\3.txt"; F = "Failed. File Already Exists!"; Style = 3 }
)

$startRow = 190
for ($i = 0; $i -lt $logRows.Count; $i++) {
    $entry = $logRows[$i]
    $r = $startRow + $i
    $target = $ws.Range("F" + $r)
    if ($entry.Style -eq 2) {
        $successTemplate.Copy($target)
    } else {
        $failTemplate.Copy($target)
    }
    # Column A holds a literal date-looking string; a leading apostrophe
    # forces it to stay plain text instead of being parsed as a date.
    $ws.Cells.Item($r, 1).Value = "'" + $entry.A
    $ws.Cells.Item($r, 2).Value = $entry.B
    $ws.Cells.Item($r, 3).Value = $entry.C
    $ws.Cells.Item($r, 4).Value = $entry.D
    $ws.Cells.Item($r, 5).Value = $entry.E
    $ws.Cells.Item($r, 6).Value = $entry.F
}

$excel.CutCopyMode = $false
